$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the number of papers/copies for the corresponding bill rows.
# These drive the existing formulas in column I (and the SUM total in I32)
# to recalculate automatically.
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1

$excel.Calculate()
